$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.001.78'
$ws.Range("E2").Value = '  +0.86%  '
$ws.Range("D3").Value = '3.180.36'
$ws.Range("E3").Value = '  -4.40%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.37'
$ws.Range("E5").Value = '  -0.44%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '172.14'
$ws.Range("E6").Value = '  -3.24%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.598'
$ws.Range("E8").Value = '  -2.58%  '
$ws.Range("D9").Value = '3.177.27'
$ws.Range("E9").Value = '  -4.44%  '
$ws.Range("E10").Value = '  -3.09%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.60'
$ws.Range("E11").Value = '  -3.97%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.394'
$ws.Range("E12").Value = '  -4.43%  '
$ws.Range("D13").Value = '3.731.18'
$ws.Range("E13").Value = '  -4.48%  '
$ws.Range("E14").Value = '  +0.90%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '27.41'
$ws.Range("E15").Value = '  -4.45%  '
$ws.Range("D16").Value = '65.947.16'
$ws.Range("E16").Value = '  +0.64%  '
$ws.Range("E17").Value = '  -2.57%  '
$ws.Range("D18").Value = '3.179.98'
$ws.Range("E18").Value = '  -4.45%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.73'
$ws.Range("E19").Value = '  +0.11%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.93'
$ws.Range("E20").Value = '  -3.62%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '361.68'
$ws.Range("E21").Value = '  -0.09%  '
$ws.Range("E22").Value = '  -1.93%  '
$ws.Range("E23").Value = '  -0.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '69.45'
$ws.Range("E24").Value = '  -2.51%  '
$ws.Range("E25").Value = '  -4.45%  '
$ws.Range("D26").Value = '3.311.64'
$ws.Range("E26").Value = '  -4.88%  '
$ws.Range("E27").Value = '  -5.89%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.87'
$ws.Range("E28").Value = '  +3.11%  '
$ws.Range("E29").Value = '  -0.15%  '
$ws.Range("E30").Value = '  -0.08%  '
$ws.Range("E31").Value = '  -0.16%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.92'
$ws.Range("E32").Value = '  -1.94%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.40'
$ws.Range("E33").Value = '  -3.91%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '22.12'
$ws.Range("E34").Value = '  -3.03%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.20'
$ws.Range("E35").Value = '  -1.36%  '
$ws.Range("E36").Value = '  -2.91%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '159.49'
$ws.Range("E37").Value = '  -0.67%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.46'
$ws.Range("E38").Value = '  -2.19%  '
$ws.Range("E39").Value = '  -0.84%  '
$ws.Range("E40").Value = '  +3.56%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '26.46'
$ws.Range("E41").Value = '  -3.39%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.49'
$ws.Range("E42").Value = '  -1.74%  '
$ws.Range("D43").Value = '2.660.49'
$ws.Range("E43").Value = '  -1.22%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '6.16'
$ws.Range("E44").Value = '  -0.95%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '4.20'
$ws.Range("E45").Value = '  -1.97%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.71'
$ws.Range("E46").Value = '  +0.22%  '
$ws.Range("E47").Value = '  -0.85%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '328.94'
$ws.Range("E48").Value = '  -1.66%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '24.05'
$ws.Range("E49").Value = '  -0.28%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0275'
$ws.Range("E50").Value = '  -1.32%  '
$ws.Range("E51").Value = '  -1.32%  '
